$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean)
$ws.Range("B3").Value = 24394.16842128429
$ws.Range("D3").Value = 1877.883891196729
$ws.Range("E3").Value = 4112.394024067734

# Row 4 (std)
$ws.Range("B4").Value = 10186.37508745993
$ws.Range("D4").Value = 1000.200061562998
$ws.Range("E4").Value = 1869.939638531676

# Row 5 (min)
$ws.Range("B5").Value = 8864.980265753431
$ws.Range("D5").Value = 307.5391671232867
$ws.Range("E5").Value = 1078.896005479453

# Row 6 (25%)
$ws.Range("B6").Value = 16292.75441095892
$ws.Range("D6").Value = 920.8718075342466
$ws.Range("E6").Value = 2147.12578630137

# Row 7 (50%)
$ws.Range("B7").Value = 22603.73099452058
$ws.Range("D7").Value = 1615.933397260274
$ws.Range("E7").Value = 4018.315223287673

# Row 8 (75%)
$ws.Range("B8").Value = 32155.58785958915
$ws.Range("D8").Value = 2783.37747671233
$ws.Range("E8").Value = 6044.991489041099

# Row 9 (max)
$ws.Range("B9").Value = 43301.83898630141
$ws.Range("D9").Value = 3509.69724383562
$ws.Range("E9").Value = 7215.295767123282

# Row 10 (Total sums)
$ws.Range("F10").Value = 35127602.52664935

# Row 11 (Residential % energy sector)
$ws.Range("G11").Value = 0.7544381176757862

# Row 12 (Community sums & % energy sector)
$ws.Range("F12").Value = 2704152.803323289
$ws.Range("G12").Value = 0.07698085291393852

# Row 13 (IGA sums & % energy sector)
$ws.Range("F13").Value = 5921847.394657535
$ws.Range("G13").Value = 0.1685810294102753
